# Regenerate the experiment task-order sheets: reorder tabs, rename them,
# and write the freshly-generated stimulus-file lists into each sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Reorder the tabs -------------------------------------------------
# Before: GNG, NB, RS, TOL, vSAT
# After : RS, GNG, vSAT, NB, TOL
# Use fresh Item() lookups by (original) name for each Move call so we
# don't operate on stale references after the collection order shifts.
$wb.Worksheets.Item("RS_TO-1651255639989003").Move($wb.Worksheets.Item("GNG_TO-16512556372807684"))
$wb.Worksheets.Item("vSAT_TO-16512556401215255").Move($wb.Worksheets.Item("NB_TO-16512556399811816"))

# --- 2) Rename the tabs to their newly generated identifiers -------------
$wb.Worksheets.Item("RS_TO-1651255639989003").Name = "RS_TO-16515890497723918"
$wb.Worksheets.Item("GNG_TO-16512556372807684").Name = "GNG_TO-1651589049803607"
$wb.Worksheets.Item("vSAT_TO-16512556401215255").Name = "vSAT_TO-1651589049866109"
$wb.Worksheets.Item("NB_TO-16512556399811816").Name = "NB_TO-16515890512542582"
$wb.Worksheets.Item("TOL_TO-16512556400449636").Name = "TOL_TO-16515890513010755"

# --- 3) Rewrite each sheet's generated stimulus-file values ---------------
# RS sheet keeps the same values ("eyes closed" / "eyes open") - no change needed.

# GNG sheet
$wsGNG = $wb.Worksheets.Item("GNG_TO-1651589049803607")
$wsGNG.Range("B2").Value = "go_stims-16515890497723918.csv"
$wsGNG.Range("B3").Value = "GNG_stims-1651589049787999.csv"
$wsGNG.Range("B4").Value = "go_stims-1651589049787999.csv"
$wsGNG.Range("B5").Value = "GNG_stims-1651589049803607.csv"

# vSAT sheet
$wsvSAT = $wb.Worksheets.Item("vSAT_TO-1651589049866109")
$wsvSAT.Range("B2").Value = "SAT_stims-16515890498192306.csv"
$wsvSAT.Range("B3").Value = "vSAT_stims-16515890498504803.csv"
$wsvSAT.Range("B4").Value = "vSAT_stims-16515890498348553.csv"
$wsvSAT.Range("B5").Value = "SAT_stims-1651589049803607.csv"

# NB sheet
$wsNB = $wb.Worksheets.Item("NB_TO-16515890512542582")
$wsNB.Range("B2").Value = "OB-16515890504950366.csv"
$wsNB.Range("B3").Value = "TB-16515890506221538.csv"
$wsNB.Range("B4").Value = "OB-16515890505909398.csv"
$wsNB.Range("B5").Value = "ZB-match_5-16515890501256156.csv"
$wsNB.Range("B6").Value = "TB-16515890508119173.csv"
$wsNB.Range("B7").Value = "ZB-match_5-16515890500787754.csv"
$wsNB.Range("B8").Value = "ZB-match_9-1651589050141241.csv"
$wsNB.Range("B9").Value = "OB-16515890505753198.csv"
$wsNB.Range("B10").Value = "TB-165158905123862.csv"

# TOL sheet
$wsTOL = $wb.Worksheets.Item("TOL_TO-16515890513010755")
$wsTOL.Range("B2").Value = "MM_stims-16515890512698603.csv"
$wsTOL.Range("B3").Value = "ZM_stims-16515890512542582.csv"
$wsTOL.Range("B4").Value = "MM_stims-16515890512854838.csv"
$wsTOL.Range("B5").Value = "ZM_stims-16515890512698603.csv"
$wsTOL.Range("B6").Value = "MM_stims-16515890513010755.csv"
$wsTOL.Range("B7").Value = "ZM_stims-16515890512854838.csv"
